$d = $word.ActiveDocument

# Find.Execute's Replacement text runs through the "typed text" autoformatter,
# which smart-quotes any straight apostrophe ('survives InsertAfter/Range.Text
# assignment unharmed, but not a Find replacement). So: do the merge replace
# with a placeholder in place of each apostrophe, then fix up the placeholder
# occurrences with a direct Range.Text assignment (no autoformat) afterwards.
function Merge-Text($searchText) {
    $placeholder = [char]1
    $placeholderText = $searchText.Replace("'", $placeholder)

    $range = $d.Content
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $placeholderText, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $searchText"
    }

    while ($true) {
        $fixRange = $d.Content
        $gotPlaceholder = $fixRange.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $gotPlaceholder) {
            break
        }
        $fixRange.Text = "'"
    }
}

# 1) Merge the three runs split by gramStart/gramEnd proofErr markers back into
#    a single run/sentence (removes the proofErr wrapper + run splits).
Merge-Text(": The dashboard should display data visualizations, such as scatter plots, line charts, and bar charts, that allow users to explore the data.")
Merge-Text("Overall, MongoDB's flexibility, scalability, and ease of use, combined with its strong support for Python, make it an excellent choice for the model component of web applications.")
Merge-Text("Anaconda: Anaconda is a distribution of the Python and R programming languages for scientific computing, that aims to simplify package management and deployment.")

# 2) Collapse the run of 12 empty, centered paragraphs before "Contact" down to
#    just 2 empty paragraphs, and remove the center alignment from those plus
#    the following Contact-block paragraphs.
$contactIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Contact") {
        $contactIndex = $i
        break
    }
}
if ($contactIndex -eq -1) {
    throw "Could not find 'Contact' paragraph"
}

# Walk backwards from the Contact paragraph collecting the contiguous empty,
# centered paragraphs that precede it (collected nearest-first / descending).
$emptyIndicesDesc = @()
$j = $contactIndex - 1
while ($j -ge 1) {
    $p = $d.Paragraphs($j)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "" -and $p.Format.Alignment -eq 1) {
        $emptyIndicesDesc += $j
        $j--
    } else {
        break
    }
}

# $emptyIndicesDesc is sorted descending (closest to Contact first). Keep the
# two paragraphs furthest from Contact (the first two in document order) and
# delete the rest, highest index first so indices stay valid while deleting.
$keepCount = 2
if ($emptyIndicesDesc.Count -gt $keepCount) {
    $toDelete = $emptyIndicesDesc[0..($emptyIndicesDesc.Count - 1 - $keepCount)]
    foreach ($idx in $toDelete) {
        $d.Paragraphs($idx).Range.Delete()
    }
}

# Re-resolve the Contact paragraph index after deletions and clear center
# alignment on it plus the two empty paragraphs before it, and on the next
# three paragraphs (name, email, university).
$contactIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Contact") {
        $contactIndex = $i
        break
    }
}
if ($contactIndex -eq -1) {
    throw "Could not find 'Contact' paragraph after delete"
}

for ($i = $contactIndex - 2; $i -le $contactIndex + 3; $i++) {
    $d.Paragraphs($i).Format.Alignment = 0
}

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
